# Apply the "changed image resize, batch size and changed loss function to
# MSE + IoU" update to the "Тесты" (Tests) log sheet: fill in the results of
# test #22 (row 23) and add a brand-new test #23 (row 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тесты")

# Test #22 (row 23): record the result of the loss-function change.
$ws.Range("H23").Value = "Train IoU: 0.51, Val IoU: 0.47. Пока что результаты от замен функции потерь незначительные. В дальнейшем нужно доработать их и собрать все влияющие факторы в одну."

# Test #23 (row 24): new row - base model swapped to EfficientNet B4, smaller batch size.
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 15
$ws.Range("F24").Value = "Замена базовой модели на EfficientNet B4, размер батча уменьшен до 16 для оптимальной скорости обучения новой модели"
$ws.Range("G24").Value = "параметры теста 4"
$ws.Range("H24").Value = "Train IoU: 0.51, Val IoU: 0.50. Результат предыдущей лучшей версии достигается на меньшем количестве эпох и размере батча, но обучение длится в разы дольше."

# Both rows now hold wrapped multi-line text, so they grow taller.
$ws.Rows.Item(23).RowHeight = 75
$ws.Rows.Item(24).RowHeight = 75

# Leave the selection on the newly-filled-in cell, like the author did.
$ws.Range("H24").Select()
